$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlCenter alignment constant, used to match style index "s=1" (horizontal
# center) applied to the pre-existing cells in this block of columns.
$xlCenter = -4108

# --- New column E ("№ поста") + revised H/I ("Количество символов" /
# "Количество изображений") stats for rows 2-27 -----------------------
$rows = @(
  @{R=2;  E=1;  H=1900; I=41},
  @{R=3;  E=2;  H=0;    I=0},
  @{R=4;  E=3;  H=15;   I=1},
  @{R=5;  E=4;  H=0;    I=0},
  @{R=6;  E=5;  H=5;    I=0},
  @{R=7;  E=6;  H=634;  I=92},
  @{R=8;  E=7;  H=1100; I=29},
  @{R=9;  E=8;  H=2600; I=401},
  @{R=10; E=9;  H=8;    I=8},
  @{R=11; E=10; H=130;  I=3},
  @{R=12; E=11; H=619;  I=23},
  @{R=13; E=12; H=50;   I=5},
  @{R=14; E=13; H=8;    I=2},
  @{R=15; E=14; H=2;    I=0},
  @{R=16; E=15; H=7800; I=435},
  @{R=17; E=16; H=12;   I=8},
  @{R=18; E=17; H=405;  I=14},
  @{R=19; E=18; H=53;   I=3},
  @{R=20; E=19; H=760;  I=28},
  @{R=21; E=20; H=8400; I=297},
  @{R=22; E=21; H=80;   I=1},
  @{R=23; E=22; H=53;   I=12},
  @{R=24; E=23; H=22;   I=0},
  @{R=25; E=24; H=118;  I=0},
  @{R=26; E=25; H=6800; I=363},
  @{R=27; E=26; H=31;   I=1}
)

foreach ($d in $rows) {
  $eCell = $ws.Cells.Item($d.R, 5)
  $eCell.Value = $d.E
  $eCell.HorizontalAlignment = $xlCenter
  $ws.Cells.Item($d.R, 8).Value = $d.H
  $ws.Cells.Item($d.R, 9).Value = $d.I
}

# --- Rows 28-31: only H/I change, column E doesn't extend this far -----
$ws.Cells.Item(28, 8).Value = 12
$ws.Cells.Item(28, 9).Value = 3

$ws.Cells.Item(29, 8).Value = 2500
$ws.Cells.Item(29, 9).Value = 88

$ws.Cells.Item(30, 8).Value = 7
$ws.Cells.Item(30, 9).Value = 2

$ws.Cells.Item(31, 8).Value = 0

# --- Row 29, columns A-D: new "Среднее" (average) row ------------------
$ws.Range("A29").Value = "Среднее"
$ws.Range("B29").Formula = "=B28/`$E`$27"
$ws.Range("C29").Formula = "=C28/`$E`$27"
$ws.Range("D29").Formula = "=D28/`$E`$27"

# --- Sheet view: scroll position + active selection ---------------------
try {
  $win = $excel.ActiveWindow
  $win.ScrollColumn = 6
  $win.ScrollRow = 1
} catch {
}
$ws.Range("D32").Select()

$wb.Application.Calculate()
